# This script re-derives 20 unique trials (by permuting the row-specific
# data H:V across the 40 data rows of the sheet) and shifts the running
# trial_total counter (column F) down by 81 for every row, matching the
# "new input files generation ... make only 20 different versions and
# duplicate many times" re-shuffle described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row number -> source row number (1-based worksheet
# row numbers). Row N after the edit receives the H:V block that used to
# live in row mapping[N] before the edit. This is a full permutation of
# rows 2..41 (some rows, like 18 and 40, map to themselves).
$mapping = @{
  2=19; 3=2; 4=33; 5=11; 6=10; 7=8; 8=34; 9=3; 10=29; 11=20; 12=13; 13=39;
  14=24; 15=17; 16=37; 17=28; 18=18; 19=16; 20=36; 21=22; 22=6; 23=30; 24=12; 25=27;
  26=25; 27=38; 28=9; 29=7; 30=35; 31=15; 32=31; 33=4; 34=14; 35=32; 36=5;
  37=41; 38=23; 39=26; 40=40; 41=21
}

$firstRow = 2
$lastRow = 41

# Read the full source block (columns H..V) once, before any writes, so
# that later writes don't clobber data that still needs to be read.
$srcRange = $ws.Range("H$firstRow`:V$lastRow")
$srcVals = $srcRange.Value()

$numRows = $lastRow - $firstRow + 1
$numCols = $srcVals.GetLength(1)

# Build the destination array applying the permutation.
$destVals = New-Object 'object[,]' $numRows, $numCols

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $mapping[$destRow]
    $destIdx = $destRow - $firstRow        # 0-based index into $destVals
    $srcIdx = $srcRow - $firstRow + 1      # 1-based index into $srcVals (COM SAFEARRAY)
    for ($c = 1; $c -le $numCols; $c++) {
        $destVals[$destIdx, ($c - 1)] = $srcVals[$srcIdx, $c]
    }
}

$destRange = $ws.Range("H$firstRow`:V$lastRow")
$destRange.Value = $destVals

# Shift the trial_total counter (column F) down by 81 for every data row;
# this column is not permuted, each row keeps its own trial number, just
# renumbered to the new, shorter running sequence.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Range("F$r")
    $cell.Value = $cell.Value() - 81
}
